$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-11 (data for 2000年-2009年), shifting the remaining rows
# (2010年-2015年, originally rows 12-17) up to rows 2-7.
$ws.Range("A2:F11").EntireRow.Delete()
